# Timesheet Calculator v2 - implement holidays
# * Mark 12-Jun-2021 (SAT) as a holiday on "Doe, Jean S." sheet:
#     - Date label gets a "(HOLIDAY)" suffix
#     - Time-in/Time-out gets populated with the standard workday slot
#     - Rendered minutes (D12) and excess/weekend minutes (F12) become 480
# * The user ends the session with the "Doe, Jean S." sheet active,
#   selection resting on A13.

$wb = $excel.ActiveWorkbook

$wsDoe = $wb.Worksheets.Item("Doe, Jean S.")

# Update the SAT (12-Jun-2021) row to reflect the holiday.
$wsDoe.Range("A12").Value = "12-Jun-2021 (HOLIDAY)"
$wsDoe.Range("C12").Value = "09:00 AM-06:00 PM"
$wsDoe.Range("D12").Value = 480
$wsDoe.Range("F12").Value = 480

# Switch the active tab from "Summary" to "Doe, Jean S." and leave the
# selection on A13, matching the state the workbook was saved in.
$wsDoe.Activate()
$wsDoe.Range("A13").Select()
